$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.652.54'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.596.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.514'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0617'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.821.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.589.49'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.523'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.635.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.275.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.843'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.786'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.734.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.909'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.46'
$ws.Range('D51').Style = 'Normal'

# --- Volume(1h) (column E) updates ---
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +1.18%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('E42').Value = '  +3.19%  '
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E46').Value = '  +8.67%  '
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('E49').Value = '  +4.54%  '
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('E51').Value = '  -0.61%  '
